$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every cell in column B that currently reads "DL-based Code Completion"
# with the new technique name "T5 model" (per commit message: "Updated DL-based with T5 model").
$lastRow = 151
for ($r = 2; $r -le $lastRow; $r++) {
    $c = $ws.Cells.Item($r, 2)
    if ($c.Text -eq "DL-based Code Completion") {
        $c.Value = "T5 model"
    }
}

# Update the active selection to match the edited range.
$ws.Range("B8:B150").Select()
